$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1 & 6. "Labelling gray images" -> "Labeling gray images"
#    Occurs twice (TOC entry + the "2.a.vi." body heading); both need
#    the identical "Labelling" -> "Labeling" substring fix, so a single
#    Replace-All handles both occurrences safely.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Labelling gray images", $true, $false, $false, $false, $false, $true, 1, $false, "Labeling gray images", 2) | Out-Null

# ---------------------------------------------------------------------
# 2 & 4. Move the lone "_GoBack" bookmark from the "2.a.ii. Resiz|ing"
#    split to inside "Visualization (Tenso|rBoard)" (first/TOC
#    occurrence only), and re-merge the now-orphaned "Resiz" + "ing"
#    runs back into a single "2.a.ii. Resizing" run.
# ---------------------------------------------------------------------
$rngTB = $d.Content
$null = $rngTB.Find.Execute("Visualization (TensorBoard)")
$tbSplit = $rngTB.Start + "Visualization (Tenso".Length
$bmRange = $d.Range($tbSplit, $tbSplit)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$rngResize = $d.Content
$null = $rngResize.Find.Execute("2.a.ii. Resiz")
$resizeStart = $rngResize.Start
$resizeFull = $d.Range($resizeStart, $resizeStart + "2.a.ii. Resizing".Length)
# Round-trip through a placeholder so the engine rebuilds a single run
# instead of leaving the old "Resiz" / "ing" run boundary in place.
$resizeFull.Text = "2.a.ii. ResizPLACEHOLDER"
$resizeFull2 = $d.Range($resizeStart, $resizeStart + "2.a.ii. ResizPLACEHOLDER".Length)
$resizeFull2.Text = "2.a.ii. Resizing"

# ---------------------------------------------------------------------
# 3. "Decod the image string from the image file." ->
#    "Decode the image string from the image file."
#    Inserted "e" becomes its own run (matches the commit's split).
# ---------------------------------------------------------------------
$rngDecod = $d.Content
$null = $rngDecod.Find.Execute("Decod the image string from the image file.")
$decodStart = $rngDecod.Start
$eInsertPoint = $d.Range($decodStart + "Decod".Length, $decodStart + "Decod".Length)
$eInsertPoint.InsertAfter("e") | Out-Null
$eChar = $d.Range($decodStart + "Decod".Length, $decodStart + "Decod".Length + 1)
$eChar.Bold = 1
$eChar.Bold = 0

# ---------------------------------------------------------------------
# 5. "... using sigmoid function ..." -> "... using Sigmoid function ..."
#    "Sigmoid" becomes its own run.
# ---------------------------------------------------------------------
$rngSig = $d.Content
$null = $rngSig.Find.Execute("Normalize the a*b* channels by using sigmoid function to evenly distribute the value of lab_ab in range (0, 1).")
$sigStart = $rngSig.Start
$sigPrefixLen = "Normalize the a*b* channels by using ".Length
$sigWordLen = "sigmoid".Length
$sigWord = $d.Range($sigStart + $sigPrefixLen, $sigStart + $sigPrefixLen + $sigWordLen)
$sigWord.Text = "Sigmoid"
$sigWord2 = $d.Range($sigStart + $sigPrefixLen, $sigStart + $sigPrefixLen + $sigWordLen)
$sigWord2.Bold = 1
$sigWord2.Bold = 0

# ---------------------------------------------------------------------
# 7. "Optimerizer: Adam" -> "Optimizer: Adam"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Optimerizer: Adam", $true, $false, $false, $false, $false, $true, 1, $false, "Optimizer: Adam", 1) | Out-Null

# ---------------------------------------------------------------------
# 8. "Comparsion between ..." -> "Comparison between ..."
#    "Comparison" becomes its own run.
# ---------------------------------------------------------------------
$rngCmp = $d.Content
$null = $rngCmp.Find.Execute("Comparsion between Basic colorization model and Full colorization model (with global network)")
$cmpStart = $rngCmp.Start
$cmpWordLen = "Comparsion".Length
$cmpWord = $d.Range($cmpStart, $cmpStart + $cmpWordLen)
$cmpWord.Text = "Comparison"
$cmpWord2 = $d.Range($cmpStart, $cmpStart + $cmpWordLen)
$cmpWord2.Bold = 1
$cmpWord2.Bold = 0

# ---------------------------------------------------------------------
# 9. "Train with both trian and test dataset" ->
#    "Train with both train and test dataset"
#    "train" becomes its own run.
# ---------------------------------------------------------------------
$rngTr = $d.Content
$null = $rngTr.Find.Execute("Train with both trian and test dataset")
$trStart = $rngTr.Start
$trPrefixLen = "Train with both ".Length
$trWordLen = "trian".Length
$trWord = $d.Range($trStart + $trPrefixLen, $trStart + $trPrefixLen + $trWordLen)
$trWord.Text = "train"
$trWord2 = $d.Range($trStart + $trPrefixLen, $trStart + $trPrefixLen + $trWordLen)
$trWord2.Bold = 1
$trWord2.Bold = 0
